# Applies the update described in the commit:
#  - Re-orders the betting-odds detail (columns F:V) for three clusters of
#    rows that represent the same set of matches but in a different order
#    (rows 89-91, rows 103-107 and rows 113-116).
#  - Appends a new match row (row 119) for Karpaty Krosno vs Unia Tarnow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: re-order a block of rows. $order[i] tells which row (relative
# to $rows) currently holds the data that should end up in $rows[i].
# We snapshot all rows first so that overlapping writes don't clobber
# data we still need to read.
# ---------------------------------------------------------------------
function Reorder-Rows {
    param($ws, [int[]]$rows, [int[]]$sourceForEachRow)

    $snapshots = @{}
    foreach ($r in $rows) {
        $snapshots[$r] = $ws.Range("F$r`:V$r").Value()
    }

    for ($i = 0; $i -lt $rows.Length; $i++) {
        $destRow = $rows[$i]
        $srcRow = $sourceForEachRow[$i]
        $ws.Range("F$destRow`:V$destRow").Value = $snapshots[$srcRow]
    }
}

# Block 1: rows 89-91 -> new row gets the data that used to sit in the
# row shown (cyclic rotation: 89<-91, 90<-89, 91<-90)
# NOTE: named parameters (-ws ... -rows ...) are not reliable in this
# PowerShell runtime, so positional arguments are used instead.
Reorder-Rows $ws @(89, 90, 91) @(91, 89, 90)

# Block 2: rows 103-107
Reorder-Rows $ws @(103, 104, 105, 106, 107) @(104, 107, 103, 105, 106)

# Block 3: rows 113-116 (simple pairwise swaps)
Reorder-Rows $ws @(113, 114, 115, 116) @(114, 113, 116, 115)

# ---------------------------------------------------------------------
# Append new row 119 (same layout/style as the existing data rows, so
# copy row 118 first to inherit formatting, then overwrite the values).
# ---------------------------------------------------------------------
$ws.Range("A118:V118").Copy($ws.Range("A119:V119"))

$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "poland"
$ws.Cells.Item(119, 3).Value = "iii-liga-group-iv"
$ws.Cells.Item(119, 4).Value = "2023-2024"
$ws.Cells.Item(119, 5).Value = 45240.70833333334
$ws.Cells.Item(119, 6).Value = "Karpaty Krosno"
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = "Unia Tarnow"
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 2.62
$ws.Cells.Item(119, 11).Value = "10/11/2023 06:12"
$ws.Cells.Item(119, 12).Value = 2.82
$ws.Cells.Item(119, 13).Value = "10/11/2023 11:12"
$ws.Cells.Item(119, 14).Value = 3.45
$ws.Cells.Item(119, 15).Value = "10/11/2023 06:12"
$ws.Cells.Item(119, 16).Value = 3.68
$ws.Cells.Item(119, 17).Value = "10/11/2023 15:03"
$ws.Cells.Item(119, 18).Value = 2.22
$ws.Cells.Item(119, 19).Value = "10/11/2023 06:12"
$ws.Cells.Item(119, 20).Value = 2.07
$ws.Cells.Item(119, 21).Value = "10/11/2023 11:12"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/poland/iii-liga-group-iv/ks-karpaty-krosno-unia-tarnow/MLmDTM5D/"
